# Update the "Förändrad" (Changed) date column (C) for data rows 2-54
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 54; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
